# "Generate Report for Handback" - localization-status report update.
#
# 1) Overview sheet: status text for both languages flips from
#    "Ready for handoff" to "Handed back: in sync with en-US", and the
#    zh-cn/de-de summary columns get wider to fit the longer text.
# 2) zh-cn / de-de detail sheets: the "Latest Target File" (I) and
#    "Latest Handback File" (J) columns get populated (they previously
#    held placeholder blanks), the "Latest Handback DateTime" (K) column
#    gets a real timestamp instead of the zero-date sentinel, and the new
#    "Latest Target File" cells become hyperlinks to the source doc
#    (mirroring the existing "Source File Name" hyperlink style/target).
#    Column widths for the file-name columns grow to fit.

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2
$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) -> cornflower blue used by the workbook's HyperLink style

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(10).ColumnWidth = 39.16

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11d3278ba4f31ef980bd1c1f3a97f14c25c93eab/e2e/a.md", "", "", "a.md")
$zhcn.Range("I2").Font.Underline = $hyperlinkUnderline
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 18:45:43"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11d3278ba4f31ef980bd1c1f3a97f14c25c93eab/e2e/a.md", "", "", "a.md")
$zhcn.Range("I3").Font.Underline = $hyperlinkUnderline
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 18:45:43"

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(10).ColumnWidth = 39.16

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11d3278ba4f31ef980bd1c1f3a97f14c25c93eab/e2e/a.md", "", "", "a.md")
$dede.Range("I2").Font.Underline = $hyperlinkUnderline
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 18:45:51"

$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/11d3278ba4f31ef980bd1c1f3a97f14c25c93eab/e2e/a.md", "", "", "a.md")
$dede.Range("I3").Font.Underline = $hyperlinkUnderline
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 18:45:51"
